# Weekly update: "Fruta / hortaliza, semanal"
# Inserts 6 new daily price rows for Ají (Vega Modelo de Temuco) ahead of the
# existing data block (old rows 770-806 shift down to 776-812), then fills
# the newly opened rows 770-775 with the new week's records.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Push the existing rows 770-806 down by 6 so rows 770-775 are free for the
# new entries (mirrors the row-shift visible across the whole diff).
$ws.Rows("770:775").Insert()

# Columns: A Mercado ID | B Mercado | C Región | D Fecha | E Codreg |
#          F Categoría ID | G Categoría | H Variedad | I Calidad | J Volumen |
#          K Precio mínimo | L Precio máximo | M Precio promedio ponderado |
#          N Unidad de comercialización | O Origen | P Precio $/Kg |
#          Q Kg o Unidades | R Clasificación
$newRows = @(
    @{ Row=770; D=44753; H="Amarillo"; I="Primera"; J=80;  K=40000; L=40000; M=40000; N="`$/caja 15 kilos"; O="Región de Arica y Parinacota"; P=2667; Q=15 },
    @{ Row=771; D=44753; H="Amarillo"; I="Segunda"; J=30;  K=30000; L=30000; M=30000; N="`$/caja 15 kilos"; O="Región de Arica y Parinacota"; P=2000; Q=15 },
    @{ Row=772; D=44753; H="Cristal";  I="Primera"; J=30;  K=33000; L=33000; M=33000; N="`$/caja 15 kilos"; O="Región de Arica y Parinacota"; P=2200; Q=15 },
    @{ Row=773; D=44753; H="Inferno";  I="Extra";   J=20;  K=22000; L=22000; M=22000; N="`$/caja 15 kilos"; O="Región de Arica y Parinacota"; P=1467; Q=15 },
    @{ Row=774; D=44753; H="Inferno";  I="Primera"; J=120; K=20000; L=20000; M=20000; N="`$/caja 15 kilos"; O="Región de Arica y Parinacota"; P=1333; Q=15 },
    @{ Row=775; D=44753; H="Inferno";  I="Segunda"; J=40;  K=15000; L=15000; M=15000; N="`$/caja 15 kilos"; O="Región de Arica y Parinacota"; P=1000; Q=15 }
)

foreach ($r in $newRows) {
    $row = $r.Row
    $ws.Cells.Item($row, 1).Value = 10
    $ws.Cells.Item($row, 2).Value = "Vega Modelo de Temuco"
    $ws.Cells.Item($row, 3).Value = "La Araucanía"
    $ws.Cells.Item($row, 4).Value = $r.D
    $ws.Cells.Item($row, 5).Value = 9
    $ws.Cells.Item($row, 6).Value = 100112021
    $ws.Cells.Item($row, 7).Value = "Ají"
    $ws.Cells.Item($row, 8).Value = $r.H
    $ws.Cells.Item($row, 9).Value = $r.I
    $ws.Cells.Item($row, 10).Value = $r.J
    $ws.Cells.Item($row, 11).Value = $r.K
    $ws.Cells.Item($row, 12).Value = $r.L
    $ws.Cells.Item($row, 13).Value = $r.M
    $ws.Cells.Item($row, 14).Value = $r.N
    $ws.Cells.Item($row, 15).Value = $r.O
    $ws.Cells.Item($row, 16).Value = $r.P
    $ws.Cells.Item($row, 17).Value = $r.Q
    $ws.Cells.Item($row, 18).Value = "Hortaliza"
}
